$wb = $excel.ActiveWorkbook

# --- Sheet references (before any renaming) ---
$sheetPI   = $wb.Worksheets.Item(1)   # "PI hours"
$sheetDept = $wb.Worksheets.Item(2)   # currently "dept hours"

# --- Step 1: duplicate the current "dept hours" sheet to the end of the
#     workbook BEFORE we overwrite its data; this duplicate will become the
#     new "unit(accumulative) hours" sheet, preserving the original
#     dept/CSL/ECE/ME/AE breakdown numbers. ---
$sheetDept.Copy($null, $sheetDept)
$sheetUnit = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetUnit.Name = "unit(accumulative) hours"

# Update the header of the new "unit(accumulative) hours" sheet.
$sheetUnit.Range("B1").Value = "unit(accumulative)"

# --- Step 2: rename the original "dept hours" sheet ---
$sheetDept.Name = "department hours"

# --- Step 3: replace "department hours" data with the new, smaller
#     department-level breakdown (ECE / ME only) ---
$sheetDept.Range("B2").Value = "ECE"
$sheetDept.Range("C2").Value = 7
$sheetDept.Range("D2").Value = 87.5

$sheetDept.Range("B3").Value = "ME"
$sheetDept.Range("C3").Value = 1
$sheetDept.Range("D3").Value = 12.5

# Remove the now-unused rows 4 and 5 (AE / extra CSL rows) so the sheet
# shrinks down to just the two data rows.
$sheetDept.Rows.Item(4).Delete()
$sheetDept.Rows.Item(4).Delete()

# --- Step 4: on "PI hours", split the old multi-valued "dept" column into
#     a single-valued "dept" column plus a new "app" column that keeps the
#     original list-style affiliations. ---

# Add new column F "app" header, copying E1's header style.
$sheetPI.Range("E1").Copy()
$sheetPI.Range("F1").PasteSpecial(-4122)
$sheetPI.Range("F1").Value = "app"

# Move the original list-valued affiliations into the new "app" column.
$sheetPI.Range("F2").Value = "['ECE', 'CSL']"
$sheetPI.Range("F3").Value = "['ME', 'AE', 'CSL']"

# Replace "dept" column values with the single primary department.
$sheetPI.Range("E2").Value = "ECE"
$sheetPI.Range("E3").Value = "ME"

$excel.CutCopyMode = $false

# Keep "PI hours" as the selected/active tab, matching the original workbook.
$sheetPI.Activate()
